$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "91.892.30"
$ws.Range("E2").Value = "  +0.79%  "

# Row 3
$ws.Range("D3").Value = "3.109.17"
$ws.Range("E3").Value = "  -0.83%  "

# Row 4
$ws.Range("E4").Value = "  -0.12%  "

# Row 5
$ws.Range("D5").Value = "'243.34"
$ws.Range("E5").Value = "  +2.71%  "

# Row 6
$ws.Range("D6").Value = "'622.52"
$ws.Range("E6").Value = "  -1.95%  "

# Row 7
$ws.Range("D7").Value = "'1.13"
$ws.Range("E7").Value = "  +5.11%  "

# Row 8
$ws.Range("D8").Value = "'0.373"
$ws.Range("E8").Value = "  +1.99%  "

# Row 9
$ws.Range("E9").Value = "  -0.09%  "

# Row 10
$ws.Range("D10").Value = "'0.765"
$ws.Range("E10").Value = "  +5.86%  "

# Row 11
$ws.Range("D11").Value = "2.728.16"
$ws.Range("E11").Value = "  -13.00%  "

# Row 12
$ws.Range("E12").Value = "  +3.18%  "

# Row 13
$ws.Range("E13").Value = "  +1.86%  "

# Row 14
$ws.Range("D14").Value = "'35.42"
$ws.Range("E14").Value = "  -2.57%  "

# Row 15
$ws.Range("D15").Value = "91.696.68"
$ws.Range("E15").Value = "  +0.75%  "

# Row 16
$ws.Range("D16").Value = "'5.47"
$ws.Range("E16").Value = "  -1.72%  "

# Row 17
$ws.Range("D17").Value = "3.700.38"

# Row 18
$ws.Range("D18").Value = "3.093.56"
$ws.Range("E18").Value = "  -1.82%  "

# Row 19
$ws.Range("E19").Value = "  -0.19%  "

# Row 20
$ws.Range("D20").Value = "'14.58"
$ws.Range("E20").Value = "  +1.66%  "

# Row 21
$ws.Range("D21").Value = "'0.0000217"
$ws.Range("E21").Value = "  +1.25%  "

# Row 22
$ws.Range("D22").Value = "'5.76"
$ws.Range("E22").Value = "  +2.10%  "

# Row 23
$ws.Range("D23").Value = "'447.24"
$ws.Range("E23").Value = "  +0.30%  "

# Row 24
$ws.Range("D24").Value = "'9.13"
$ws.Range("E24").Value = "  +1.64%  "

# Row 25
$ws.Range("E25").Value = "  -1.48%  "

# Row 26
$ws.Range("D26").Value = "'90.63"
$ws.Range("E26").Value = "  +0.47%  "

# Row 27
$ws.Range("D27").Value = "'11.94"
$ws.Range("E27").Value = "  -3.78%  "

# Row 30
$ws.Range("D30").Value = "'0.183"
$ws.Range("E30").Value = "  +14.28%  "

# Row 31
$ws.Range("D31").Value = "'0.236"
$ws.Range("E31").Value = "  +18.32%  "

# Row 32
$ws.Range("D32").Value = "'9.38"
$ws.Range("E32").Value = "  -3.30%  "

# Row 33
$ws.Range("E33").Value = "  +1.36%  "

# Row 34
$ws.Range("E34").Value = "  +12.91%  "

# Row 35
$ws.Range("E35").Value = "  +31.50%  "

# Row 36
$ws.Range("D36").Value = "'26.60"
$ws.Range("E36").Value = "  -1.15%  "

# Row 37
$ws.Range("D37").Value = "'7.63"
$ws.Range("E37").Value = "  +6.74%  "

# Row 38
$ws.Range("D38").Value = "'4.11"
$ws.Range("E38").Value = "  +21.37%  "

# Row 39
$ws.Range("E39").Value = "  -0.92%  "

# Row 40
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").Value = "'3.63"
$ws.Range("E40").Value = "  -4.60%  "

# Row 41
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").Value = "'489.45"
$ws.Range("E41").Value = "  -4.82%  "

# Row 42
$ws.Range("D42").Value = "'1.29"
$ws.Range("E42").Value = "  -1.27%  "

# Row 43
$ws.Range("E43").Value = "  +0.60%  "

# Row 44
$ws.Range("D44").Value = "'22.15"
$ws.Range("E44").Value = "  -0.26%  "

# Row 45
$ws.Range("E45").Value = "  +0.01%  "

# Row 46
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").Value = "'1.91"
$ws.Range("E46").Value = "  -1.55%  "

# Row 47
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").Value = "'154.48"
$ws.Range("E47").Value = "  +2.27%  "

# Row 48
$ws.Range("E48").Value = "  -0.95%  "

# Row 49
$ws.Range("E49").Value = "  -0.42%  "

# Row 50
$ws.Range("E50").Value = "  -0.57%  "

# Row 51
$ws.Range("D51").Value = "'44.62"
$ws.Range("E51").Value = "  -2.59%  "
